# Mifos loan product workbook: fix product short-name string (insert a
# dash after "296") on both the input and output sheets, and leave the
# output sheet ("ProductLoanOutput") as the active tab/selection.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$productCode = "296-MS-EPP-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-ONTIME"

$ws1.Range("B1").Value = $productCode
$ws2.Range("B1").Value = $productCode

# Move the selection/active sheet: ProductLoanInput's selection moves to
# B1 (no longer the selected tab), and ProductLoanOutput becomes the
# active tab with its selection on B1.
$ws1.Activate()
$ws1.Range("B1").Select()

$ws2.Activate()
$ws2.Range("B1").Select()
